$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh: insert a new data row at row 6 (shifting the existing
# historical rows 6-20 down to rows 7-21) and populate it with the new
# week's record.
$ws.Rows.Item(6).Insert()

$ws.Range("A6").Value = 7
$ws.Range("B6").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C6").Value = "Ñuble"
$ws.Range("D6").Value = 45274
$ws.Range("E6").Value = 16
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100103
$ws.Range("H6").Value = "Frutos de hueso (carozo)"
$ws.Range("I6").Value = 100103003
$ws.Range("J6").Value = "Damasco"
$ws.Range("K6").Value = "Castle Brite"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 50
$ws.Range("N6").Value = 20000
$ws.Range("O6").Value = 20000
$ws.Range("P6").Value = 20000
$ws.Range("Q6").Value = "$/bandeja 10 kilos"
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 2000
$ws.Range("T6").Value = 10
